$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the FAILED/empty/date cells from row 2 (E2:G2)
$ws.Range("E2:G2").ClearContents()

# Update the selection to H6
$ws.Range("H6").Select()
